$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining cells of row 7 (PriceChange, UpDown)
$ws.Range("X7").Value = -3.5699769999999944
$ws.Range("Y7").Value = "Down"

# Add new row 8 with the full data set
$ws.Range("A8").Value = 42649.8909375
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 40
$ws.Range("E8").Value = 31470
$ws.Range("F8").Value = 3785
$ws.Range("G8").Value = 59
$ws.Range("H8").Value = 39
$ws.Range("I8").Value = 84
$ws.Range("J8").Value = 14
$ws.Range("K8").Value = 29646
$ws.Range("L8").Value = 457
$ws.Range("M8").Value = 301
$ws.Range("N8").Value = 175
$ws.Range("O8").Value = 30
$ws.Range("P8").Value = "Bag"
$ws.Range("Q8").Value = 49.72799223503381
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0.1095
$ws.Range("S8").NumberFormat = $ws.Range("S7").NumberFormat
$ws.Range("T8").Value = 0.0025000000000000001
$ws.Range("T8").NumberFormat = $ws.Range("T7").NumberFormat
$ws.Range("U8").Value = 5.95
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0
